$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Consequence of “z”").Name = "“z” consequence"
$wb.Worksheets.Item("Consequence of “x”").Name = "“x” consequence"
$wb.Worksheets.Item("Consequence of “c”").Name = "“c” consequence"
$wb.Worksheets.Item("Consequence of “v”").Name = "“v” consequence"
$wb.Worksheets.Item("Consequence of “b”").Name = "“b” consequence"
